$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite all cell values (rows 1-33) per the new Tu Vi interpretation content
$ws.Range("A1").Value = "Tử Vi"
$ws.Range("A2").Value = "Vô Chính Diệu"
$ws.Range("B2").Value = "Bạn là người khôn ngoan, sắc sảo , nếu là con cả thì sức khoẻ thường kém, lúc thiếu thời vất vả."
$ws.Range("C2").Value = "Hay phải phiêu bạt, chỗ ở không cố định."
$ws.Range("D2").Value = "Tâm lý của bạn hay bị ảnh hưởng bên ngoài nhưng thực chất đó cũng chỉ là do bạn muốn tiếp thu ý kiến mọi người."
$ws.Range("A3").Value = "Mệnh Không Thân Kiếp"
$ws.Range("B3").Value = "Bạn là người khôn ngoan, sắc sảo nên ông trời thử thách bạn với những hoàn cảnh trớ trêu."
$ws.Range("C3").Value = "Bạn cũng là tuýp người nhiệt tình 5 phút, cả thèm mau chán."
$ws.Range("A4").Value = "Tử Phủ Vũ Tướng"
$ws.Range("B4").Value = "Bạn là người có tính cách thích lãnh đạo, chỉ đạo."
$ws.Range("A5").Value = "Sát Phá Tham hội chiếu tại mệnh"
$ws.Range("B5").Value = "Bạn là người thiên về kinh doanh thương mại. "
$ws.Range("C5").Value = "Bản tính có xu hướng sát phạt, hơn thua, thích thay cũ đổi mới và có nhiều ham muốn."
$ws.Range("A6").Value = "Sát Phá Tham"
$ws.Range("A7").Value = "Tham Hỏa Linh"
$ws.Range("A8").Value = "Cơ Nguyệt Đồng Lương"
$ws.Range("B8").Value = "Bạn là tuýp người nhẹ nhàng, có nội tâm sâu sắc, phù hợp với môi trường giáo dục, công việc cần chuyên lý thuyết và tư duy cao."
$ws.Range("A9").Value = "Quang Quý"
$ws.Range("B9").Value = "Bạn sẽ gặp nhiều may mắn"
$ws.Range("A10").Value = "Tang Hổ hội chiếu tại Mệnh"
$ws.Range("B10").Value = "Bạn là người hay lo lắng, suy nghĩ, cuộc sống vất vả nên đôi khi sẽ lo lắng quá thực tế."
$ws.Range("C10").Value = "Bạn cũng là tuýp người có trách nhiệm với bản thân và người xung quanh."
$ws.Range("A11").Value = "Bạch Hổ tọa thủ tại Mệnh"
$ws.Range("B11").Value = "Cuộc sống của bạn xuất thân nghèo khó"
$ws.Range("C11").Value = "Bạn là người hay suy nghĩ lo lắng, tuy nhiên sự lo lắng của bạn là có cơ sở."
$ws.Range("D11").Value = "Bạn là người can đảm, cương nghị, ương ngạch."
$ws.Range("A12").Value = "Vũ Khúc"
$ws.Range("A13").Value = "Tham Lang"
$ws.Range("A14").Value = "Vũ Tham"
$ws.Range("A15").Value = "Vũ Khúc và Tham Lang đồng cung tại Mùi"
$ws.Range("B15").Value = "Cuộc sống trước nghèo sau giàu. Hơn 30 tuổi trở đi mới khá giả"
$ws.Range("A16").Value = "Phong Cáo"
$ws.Range("B16").Value = "Bạn là người có bằng cấp, đỗ đạt cao."
$ws.Range("A17").Value = "Mệnh Tý Ngọ có Thiên Khốc Thiên Hư đồng cung"
$ws.Range("B17").Value = "Thiếu thời nghèo túng, trung niên khá giả, về già giàu có."
$ws.Range("A18").Value = "Hỏa Linh"
$ws.Range("B18").Value = "Tính tình nóng giận, liều lĩnh"
$ws.Range("A19").Value = "Linh Tinh"
$ws.Range("B19").Value = "Tính tình nóng giận"
$ws.Range("A20").Value = "Hỏa Tinh hội chiếu tại mệnh"
$ws.Range("B20").Value = "Tính tình nóng giận"
$ws.Range("A21").Value = "Linh Hỏa"
$ws.Range("B21").Value = "Tính tình nóng giận, liều lĩnh, bướng"
$ws.Range("A22").Value = "Tử Vi tọa thủ cung Mệnh ở Tỵ, Ngọ, Dần, Thân"
$ws.Range("B22").Value = "Bạn là người Thông minh, trung hậu."
$ws.Range("A23").Value = "Tử Vi tọa thủ cung Mệnh ở Thìn, Tuất"
$ws.Range("B23").Value = "Bạn là người đa mưu, túc trí nhưng vì cái lợi bản thân là phần nhiều."
$ws.Range("A24").Value = "Tử Vi tọa thủ cung Mệnh ở Sửu, Mùi"
$ws.Range("B24").Value = "Bạn là người thông minh, mưu lược, nhưng có phần liều lĩnh."
$ws.Range("A25").Value = "Tử Vi tọa thủ cung Mệnh ở Tý, Hợi, Mão, Dậu"
$ws.Range("B25").Value = "Bạn hơi kém thông minh, nhưng bản tính đôn hậu."
$ws.Range("C25").Value = "Quyền uy kém rực rỡ, khả năng tiêu giảm tai ách bị giảm nhiều."
$ws.Range("A26").Value = "Tử Vi tọa thủ cung Mệnh và hội chiếu các sao Thiên Tướng, Văn Khúc, Văn Xương, Thiên Khôi, Thiên Việt, Tả Phù, Hữu Bật"
$ws.Range("B26").Value = "Bạn là người có uy quyền khiến người khác nể trọng và giúp đỡ. Bản thân ra ngoài gặp nhiều may mắn."
$ws.Range("A27").Value = "Tử Vi tọa thủ cung Mệnh và hội chiếu Thiên Phủ"
$ws.Range("B27").Value = "Bạn có nhiều tiền bạc, của cải."
$ws.Range("A28").Value = "Tử Vi tọa thủ cung Mệnh và gặp Thiên Mã, Lộc Tồn"
$ws.Range("B28").Value = "Độ số quyền lực của bạn được tăng thêm."
$ws.Range("A29").Value = "Tử Vi đồng cung với Thất Sát"
$ws.Range("B29").Value = "Độ số quyền lực của bạn là tuyệt đối."
$ws.Range("C29").Value = "Chế ác được sự tác họa của Hỏa Linh"
$ws.Range("A30").Value = "Tử Vi tọa thủ cung Mệnh và gặp Kình Dương, Đà La"
$ws.Range("B30").Value = "Bạn như vị vua bị vậy hãm."
$ws.Range("C30").Value = "Bị tiểu nhân làm hại."
$ws.Range("A31").Value = "Tử Vi tọa thủ cung Mệnh và gặp Địa Không, Địa Kiếp"
$ws.Range("B31").Value = "Bạn như vị vua bị vây hãm."
$ws.Range("C31").Value = "Bị tiểu nhân làm hại."
$ws.Range("A32").Value = "Tử Vi tọa thủ cung Mệnh gặp Kình Dương hoặc Đà La hoặc Địa Không hoặc Địa Kiếp"
$ws.Range("B32").Value = "Bạn như vị vua bị vây hãm."
$ws.Range("C32").Value = "Bị tiểu nhân làm hại."
$ws.Range("A33").Value = "Tử Vi tọa thủ cung Mệnh và gặp Kình Dương, Đà La, Địa Không, Địa Kiếp"
$ws.Range("B33").Value = "Bạn như vị vua bị vây hãm không lối thoát."
$ws.Range("C33").Value = "Bị tiểu nhân làm hại."

# Highlight the new "cach cuc" (pattern) rows A22:A33 with a yellow fill
$ws.Range("A22:A33").Interior.Color = 65535

# Widen column A to fit the new, longer pattern descriptions (target stored
# width 32.42578125 chars; 31.65 is the closest input this engine's column
# width quantizer resolves to that value)
$ws.Columns("A").ColumnWidth = 31.65

# Restore the user's last on-screen selection
$ws.Range("G29").Select()
